$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1169.6428
$ws.Range("I19").Value = 2118.5833
$ws.Range("J19").Value = 457.9375
$ws.Range("K19").Value = 2118.5833
$ws.Range("L19").Value = 457.9375
$ws.Range("M19").Value = -1943.5833
$ws.Range("N19").Value = -807.9375
$ws.Range("H40").Value = 1233.3334
$ws.Range("I40").Value = 1100
$ws.Range("K40").Value = 1100
$ws.Range("M40").Value = -925
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H97").Value = 538.5
$ws.Range("J97").Value = 538.5
$ws.Range("L97").Value = 1615.5
$ws.Range("N97").Value = -2607.5
$ws.Range("H112").Value = 1342.9166
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1491.5
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 4474.5
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -6690.5
$ws.Range("H137").Value = 50002644
$ws.Range("I137").Value = 1900
$ws.Range("K137").Value = 5700
$ws.Range("M137").Value = -3150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1344.3529
$ws.Range("I45").Value = 908.4286
$ws.Range("J45").Value = 1649.5
$ws.Range("K45").Value = 908.4286
$ws.Range("L45").Value = 1649.5
$ws.Range("M45").Value = -531.4286
$ws.Range("N45").Value = -2403.5
$ws.Range("H74").Value = 2665.1
$ws.Range("I74").Value = 970.625
$ws.Range("J74").Value = 4601.643
$ws.Range("K74").Value = 970.625
$ws.Range("L74").Value = 4601.643
$ws.Range("M74").Value = -96.625
$ws.Range("N74").Value = -6349.643
$ws.Range("H77").Value = 2665.1
$ws.Range("I77").Value = 970.625
$ws.Range("J77").Value = 4601.643
$ws.Range("K77").Value = 4853.125
$ws.Range("L77").Value = 23008.215
$ws.Range("M77").Value = -485.125
$ws.Range("N77").Value = -31744.215
$ws.Range("H132").Value = 3220.8965
$ws.Range("I132").Value = 2973.5833
$ws.Range("K132").Value = 8920.749899999999
$ws.Range("M132").Value = -6390.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7836.1177
$ws.Range("I134").Value = 8575.066000000001
$ws.Range("K134").Value = 25725.198
$ws.Range("M134").Value = -23190.198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2189.9473
$ws.Range("I31").Value = 1678.1
$ws.Range("J31").Value = 2758.6667
$ws.Range("K31").Value = 1678.1
$ws.Range("L31").Value = 2758.6667
$ws.Range("M31").Value = -1383.1
$ws.Range("N31").Value = -3348.6667
$ws.Range("H34").Value = 2189.9473
$ws.Range("I34").Value = 1678.1
$ws.Range("J34").Value = 2758.6667
$ws.Range("K34").Value = 1678.1
$ws.Range("L34").Value = 2758.6667
$ws.Range("M34").Value = -1476.1
$ws.Range("N34").Value = -3162.6667
$ws.Range("H86").Value = 8002.8
$ws.Range("I86").Value = 8507
$ws.Range("J86").Value = 7666.6665
$ws.Range("K86").Value = 8507
$ws.Range("L86").Value = 7666.6665
$ws.Range("M86").Value = -7384
$ws.Range("N86").Value = -9912.666499999999
$ws.Range("H89").Value = 8002.8
$ws.Range("I89").Value = 8507
$ws.Range("J89").Value = 7666.6665
$ws.Range("K89").Value = 42535
$ws.Range("L89").Value = 38333.3325
$ws.Range("M89").Value = -36919
$ws.Range("N89").Value = -49565.3325
$ws.Range("H99").Value = 81056.31
$ws.Range("I99").Value = 37848.43
$ws.Range("J99").Value = 202038.4
$ws.Range("K99").Value = 37848.43
$ws.Range("L99").Value = 202038.4
$ws.Range("M99").Value = -36350.43
$ws.Range("N99").Value = -205034.4
$ws.Range("H126").Value = 81056.31
$ws.Range("I126").Value = 37848.43
$ws.Range("J126").Value = 202038.4
$ws.Range("K126").Value = 113545.29
$ws.Range("L126").Value = 606115.2
$ws.Range("M126").Value = -111075.29
$ws.Range("N126").Value = -611055.2
$ws.Range("H134").Value = 2662.923
$ws.Range("I134").Value = 2497.4707
$ws.Range("K134").Value = 7492.4121
$ws.Range("M134").Value = -4957.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 844.56665
$ws.Range("I5").Value = 556.2632
$ws.Range("J5").Value = 978.1707
$ws.Range("K5").Value = 1668.7896
$ws.Range("L5").Value = 2934.5121
$ws.Range("M5").Value = -1556.7896
$ws.Range("N5").Value = -3158.5121
$ws.Range("H40").Value = 206.375
$ws.Range("I40").Value = 212.75
$ws.Range("K40").Value = 851
$ws.Range("M40").Value = -782
$ws.Range("H69").Value = 969.4167
$ws.Range("I69").Value = 760.2
$ws.Range("J69").Value = 1118.8572
$ws.Range("K69").Value = 2280.6
$ws.Range("L69").Value = 3356.5716
$ws.Range("M69").Value = -1469.6
$ws.Range("N69").Value = -4978.571599999999
$ws.Range("H72").Value = 969.4167
$ws.Range("I72").Value = 760.2
$ws.Range("J72").Value = 1118.8572
$ws.Range("K72").Value = 6841.8
$ws.Range("L72").Value = 10069.7148
$ws.Range("M72").Value = -2785.8
$ws.Range("N72").Value = -18181.7148
$ws.Range("H131").Value = 2716.862
$ws.Range("J131").Value = 1741.9286
$ws.Range("L131").Value = 5225.7858
$ws.Range("N131").Value = -15305.7858
$ws.Range("H135").Value = 844.56665
$ws.Range("I135").Value = 556.2632
$ws.Range("J135").Value = 978.1707
$ws.Range("K135").Value = 5006.3688
$ws.Range("L135").Value = 8803.5363
$ws.Range("M135").Value = -2471.3688
$ws.Range("N135").Value = -13873.5363

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
$ws.Range("H107").Value = 410.56668
$ws.Range("I107").Value = 361.38095
$ws.Range("J107").Value = 525.3333
$ws.Range("K107").Value = 361.38095
$ws.Range("L107").Value = 525.3333
$ws.Range("M107").Value = 1558.61905
$ws.Range("N107").Value = -4365.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 905
$ws.Range("I46").Value = 903.2258
$ws.Range("J46").Value = 912.8570999999999
$ws.Range("K46").Value = 903.2258
$ws.Range("L46").Value = 912.8570999999999
$ws.Range("M46").Value = -715.2258
$ws.Range("N46").Value = -1288.8571
$ws.Range("H61").Value = 1606.25
$ws.Range("I61").Value = 1490
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 1490
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -1288
$ws.Range("N61").Value = -2204
$ws.Range("H93").Value = 1788.0555
$ws.Range("I93").Value = 1484.0714
$ws.Range("J93").Value = 2852
$ws.Range("K93").Value = 1484.0714
$ws.Range("L93").Value = 2852
$ws.Range("M93").Value = -236.0714
$ws.Range("N93").Value = -5348
$ws.Range("H113").Value = 1606.25
$ws.Range("I113").Value = 1490
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1490
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 680
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 2140.5454
$ws.Range("I122").Value = 1699.5
$ws.Range("J122").Value = 2669.8
$ws.Range("K122").Value = 5098.5
$ws.Range("L122").Value = 8009.400000000001
$ws.Range("M122").Value = -2648.5
$ws.Range("N122").Value = -12909.4
$ws.Range("H132").Value = 5275.5312
$ws.Range("I132").Value = 6704.7144
$ws.Range("J132").Value = 2547.0908
$ws.Range("K132").Value = 20114.1432
$ws.Range("L132").Value = 7641.2724
$ws.Range("M132").Value = -17584.1432
$ws.Range("N132").Value = -12701.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2107.4856
$ws.Range("I122").Value = 1596.6957
$ws.Range("J122").Value = 3086.5
$ws.Range("K122").Value = 4790.0871
$ws.Range("L122").Value = 9259.5
$ws.Range("M122").Value = -2340.0871
$ws.Range("N122").Value = -14159.5
$ws.Range("H132").Value = 2564.8914
$ws.Range("I132").Value = 3038.138
$ws.Range("J132").Value = 1757.5883
$ws.Range("K132").Value = 9114.414000000001
$ws.Range("L132").Value = 5272.7649
$ws.Range("M132").Value = -6584.414000000001
$ws.Range("N132").Value = -10332.7649
